$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '27.825.35'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.628.44'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.21'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.26%  '
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0608'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '1.656.83'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.556'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').Value = '27.839.38'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').Value = '0.0₃0718'
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.996'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.27%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0480'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.11'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').Value = '1.412.52'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.850'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.11%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('D46').Value = '1.767.19'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('E47').Value = '  -3.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('E49').Value = '  +0.84%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('E51').Value = '  +0.81%  '
